# Revert merge: restore C7 (password column, row for uid S1234567A) back to
# "NEWPASSWORD" (undoing the value that had been merged in as "Password").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = "NEWPASSWORD"
